$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("GSE")

# The old column D ("Terms Typically Offered") needs to become column G.
# Insert three new blank columns ahead of it for Corequisites / Concurrent / Recommended;
# this shifts the existing "Terms Typically Offered" data from D to G automatically.
$ws.Range("D1:F1").EntireColumn.Insert()

# Header row
$ws.Range("C1").Value = "Prerequisites"
$ws.Range("D1").Value = "Corequisites"
$ws.Range("E1").Value = "Concurrent"
$ws.Range("F1").Value = "Recommended"
$ws.Range("G1").Value = "Terms Typically Offered"

# GSE 500
$ws.Range("C2").Value = "Consent of department head."
$ws.Range("D2").Value = "NA"
$ws.Range("E2").Value = "NA"
$ws.Range("F2").Value = "NA"
$ws.Range("G2").Value = "TBD"

# GSE 510
$ws.Range("C3").Value = "Graduate standing."
$ws.Range("D3").Value = "NA"
$ws.Range("E3").Value = "NA"
$ws.Range("F3").Value = "MATH 206, MATH 244, or MATH 408."
$ws.Range("G3").Value = "SU "

# GSE 511
$ws.Range("C4").Value = "Concurrent enrollment in GSE 510 and graduate standing."
$ws.Range("D4").Value = "NA"
$ws.Range("E4").Value = "NA"
$ws.Range("F4").Value = "NA"
$ws.Range("G4").Value = "F"

# GSE 512
$ws.Range("C5").Value = "GSE 511 and graduate standing."
$ws.Range("D5").Value = "NA"
$ws.Range("E5").Value = "NA"
$ws.Range("F5").Value = "NA"
$ws.Range("G5").Value = "W"

# GSE 518
$ws.Range("C6").Value = "Graduate standing."
$ws.Range("D6").Value = "NA"
$ws.Range("E6").Value = "NA"
$ws.Range("F6").Value = "MATH 206 or MATH 244 or GSE 510."
$ws.Range("G6").Value = "SU "

# GSE 520
$ws.Range("C7").Value = "GSE 518 and graduate standing."
$ws.Range("D7").Value = "NA"
$ws.Range("E7").Value = "NA"
$ws.Range("F7").Value = "ECON 339."
$ws.Range("G7").Value = "F "

# GSE 522
$ws.Range("C8").Value = "GSE 520 and graduate standing."
$ws.Range("D8").Value = "NA"
$ws.Range("E8").Value = "NA"
$ws.Range("F8").Value = "NA"
$ws.Range("G8").Value = "W"

# GSE 524
$ws.Range("C9").Value = "Graduate standing."
$ws.Range("D9").Value = "NA"
$ws.Range("E9").Value = "NA"
$ws.Range("F9").Value = "NA"
$ws.Range("G9").Value = "F"

# GSE 526
$ws.Range("C10").Value = "GSE 520 and graduate standing."
$ws.Range("D10").Value = "NA"
$ws.Range("E10").Value = "NA"
$ws.Range("F10").Value = "GSE 524."
$ws.Range("G10").Value = "W "

# GSE 532
$ws.Range("C11").Value = "GSE 511 and graduate standing."
$ws.Range("D11").Value = "NA"
$ws.Range("E11").Value = "NA"
$ws.Range("F11").Value = "NA"
$ws.Range("G11").Value = "TBD"

# GSE 534
$ws.Range("C12").Value = "GSE 511 and graduate standing."
$ws.Range("D12").Value = "NA"
$ws.Range("E12").Value = "NA"
$ws.Range("F12").Value = "NA"
$ws.Range("G12").Value = "TBD"

# GSE 536
$ws.Range("C13").Value = "GSE 511 and graduate standing."
$ws.Range("D13").Value = "NA"
$ws.Range("E13").Value = "NA"
$ws.Range("F13").Value = "NA"
$ws.Range("G13").Value = "TBD"

# GSE 538
$ws.Range("C14").Value = "GSE 511 and graduate standing."
$ws.Range("D14").Value = "NA"
$ws.Range("E14").Value = "NA"
$ws.Range("F14").Value = "NA"
$ws.Range("G14").Value = "SP"

# GSE 542
$ws.Range("C15").Value = "GSE 522 and graduate standing, or consent of instructor."
$ws.Range("D15").Value = "NA"
$ws.Range("E15").Value = "NA"
$ws.Range("F15").Value = "GSE 526."
$ws.Range("G15").Value = "SP "

# GSE 544
$ws.Range("C16").Value = "GSE 520 and graduate standing."
$ws.Range("D16").Value = "NA"
$ws.Range("E16").Value = "NA"
$ws.Range("F16").Value = "GSE 524 and GSE 526."
$ws.Range("G16").Value = "SP "

# GSE 570
$ws.Range("C17").Value = "Graduate standing or consent of instructor."
$ws.Range("D17").Value = "NA"
$ws.Range("E17").Value = "NA"
$ws.Range("F17").Value = "NA"
$ws.Range("G17").Value = "TBD"

# GSE 580
$ws.Range("C18").Value = "Graduate standing."
$ws.Range("D18").Value = "NA"
$ws.Range("E18").Value = "NA"
$ws.Range("F18").Value = "NA"
$ws.Range("G18").Value = "TBD"

# GSE 599
$ws.Range("C19").Value = "Graduate standing and consent of thesis committee."
$ws.Range("D19").Value = "NA"
$ws.Range("E19").Value = "NA"
$ws.Range("F19").Value = "NA"
$ws.Range("G19").Value = "TBD"
